$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 209 (shifts existing rows 209-243 down to 210-244)
$ws.Rows(209).Insert()

# Populate the new row 209 with the new weekly record
$ws.Cells.Item(209, 1).Value = 4
$ws.Cells.Item(209, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(209, 3).Value = "Los Lagos"
$ws.Cells.Item(209, 4).Value = 44694
$ws.Cells.Item(209, 5).Value = 10
$ws.Cells.Item(209, 6).Value = 100112044
$ws.Cells.Item(209, 7).Value = "Perejil"
$ws.Cells.Item(209, 8).Value = "Sin especificar"
$ws.Cells.Item(209, 9).Value = "Primera"
$ws.Cells.Item(209, 10).Value = 160
$ws.Cells.Item(209, 11).Value = 5000
$ws.Cells.Item(209, 12).Value = 5000
$ws.Cells.Item(209, 13).Value = 5000
$ws.Cells.Item(209, 14).Value = "`$/docena de atados (3 kilos)"
$ws.Cells.Item(209, 15).Value = "Región Metropolitana"
$ws.Cells.Item(209, 16).Value = 1667
$ws.Cells.Item(209, 17).Value = 3
$ws.Cells.Item(209, 18).Value = "Hortaliza"
